$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.803.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.889.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3801"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07711"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9621"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.896.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.968"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.677"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07066"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "83.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009501"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.782.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.365"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.147.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.637"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.815"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09262"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8528"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.073"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.245"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.074"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.156"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05643"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.006"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02039"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5510"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.427"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1752"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000002879"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -25.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.250"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.686"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5176"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.060"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06760"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.49%  "
